# Add the new "Input" sheet after the last existing sheet (so it becomes
# the 3rd sheet, matching sheetId="3" / r:id="rId3" in the target workbook).
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Input"

# Populate the new sheet: A1 = "Carrier" (reuses existing shared string),
# A2 = "UPS" (new shared string).
$newSheet.Range("A1").Value = "Carrier"
$newSheet.Range("A2").Value = "UPS"

# A1 gets a text number format plus a thin box border around it.
$headerCell = $newSheet.Range("A1")
$headerCell.NumberFormat = "@"
$headerCell.Borders.Weight = 2

# Leave the new sheet's selection on B2 (matches the recorded selection).
$newSheet.Range("B2").Select()

# Switch back to Sheet1 as the active tab, with A2 selected.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()
$sheet1.Range("A2").Select()
